$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (closest values achievable through the ColumnWidth
# pixel-grid rounding that maps onto the target stored widths)
$ws.Columns.Item(6).ColumnWidth = 2.334635416666667
$ws.Columns.Item(9).ColumnWidth = 1.334635416666667
$ws.Columns.Item(11).ColumnWidth = 6.834635416666667
$ws.Columns.Item(14).ColumnWidth = 3.834635416666667

# Update row 1 values
$ws.Range("C1").Value = 19
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 32
$ws.Range("F1").Value = 13
$ws.Range("G1").Value = 14
$ws.Range("H1").Value = 30
$ws.Range("I1").Value = 3
$ws.Range("J1").Value = 15
$ws.Range("K1").Value = 0.082009999999999986
$ws.Range("L1").Value = 0.051000000000000004
$ws.Range("M1").Value = 0.032000000000000001
$ws.Range("N1").Value = 0.080000000000000002
